$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Every data row (2..482) had its "Förändrad" (column C) date bumped
#    from 45202 to 45203 (one day later).
for ($r = 2; $r -le 482; $r++) {
    $cur = $ws.Cells.Item($r, 3).Value2()
    $ws.Cells.Item($r, 3).Value = $cur + 1
}

# 2. Rows 3 and 4 swapped places (A 31224-2022 now comes before
#    A 22996-2019), and the A 31224-2022 record picked up extra
#    species / updated counts. Rewrite both rows fully to reflect the
#    new ordering and values.

$nl = "`r`n"

# New row 3: A 31224-2022 (updated counts & species list)
$ws.Cells.Item(3, 1).Value = "A 31224-2022"
$ws.Cells.Item(3, 2).Value = 44771
$ws.Cells.Item(3, 3).Value = 45203
$ws.Cells.Item(3, 4).Value = "GÄVLEBORGS LÄN"
$ws.Cells.Item(3, 5).Value = "NORDANSTIG"
$ws.Cells.Item(3, 7).Value = 1.2
$ws.Cells.Item(3, 8).Value = 1
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 3
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 6
$ws.Cells.Item(3, 16).Value = 3
$ws.Cells.Item(3, 17).Value = 7
$ws.Cells.Item(3, 18).Value = "Bitter taggsvamp" + $nl + "Knärot" + $nl + "Lammticka" + $nl + "Orange taggsvamp" + $nl + "Svart taggsvamp" + $nl + "Talltaggsvamp" + $nl + "Tallfingersvamp"
$ws.Cells.Item(3, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/artfynd/A 31224-2022.xlsx", "A 31224-2022")'
$ws.Cells.Item(3, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/kartor/A 31224-2022.png", "A 31224-2022")'
$ws.Cells.Item(3, 21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/knärot/A 31224-2022.png", "A 31224-2022")'
$ws.Cells.Item(3, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/klagomål/A 31224-2022.docx", "A 31224-2022")'
$ws.Cells.Item(3, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/klagomålsmail/A 31224-2022.docx", "A 31224-2022")'
$ws.Cells.Item(3, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/tillsyn/A 31224-2022.docx", "A 31224-2022")'
$ws.Cells.Item(3, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/tillsynsmail/A 31224-2022.docx", "A 31224-2022")'

# New row 4: A 22996-2019 (unchanged values, just moved down one row)
$ws.Cells.Item(4, 1).Value = "A 22996-2019"
$ws.Cells.Item(4, 2).Value = 43591
$ws.Cells.Item(4, 3).Value = 45203
$ws.Cells.Item(4, 4).Value = "GÄVLEBORGS LÄN"
$ws.Cells.Item(4, 5).Value = "NORDANSTIG"
$ws.Cells.Item(4, 7).Value = 11.1
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 3
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 2
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 5
$ws.Cells.Item(4, 18).Value = "Lunglav" + $nl + "Violettgrå tagellav" + $nl + "Korallblylav" + $nl + "Sotriska" + $nl + "Stuplav"
$ws.Cells.Item(4, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/artfynd/A 22996-2019.xlsx", "A 22996-2019")'
$ws.Cells.Item(4, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/kartor/A 22996-2019.png", "A 22996-2019")'
$ws.Cells.Item(4, 21).ClearContents()
$ws.Cells.Item(4, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/klagomål/A 22996-2019.docx", "A 22996-2019")'
$ws.Cells.Item(4, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/klagomålsmail/A 22996-2019.docx", "A 22996-2019")'
$ws.Cells.Item(4, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/tillsyn/A 22996-2019.docx", "A 22996-2019")'
$ws.Cells.Item(4, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_NORDANSTIG/tillsynsmail/A 22996-2019.docx", "A 22996-2019")'

# Row heights stay at their original explicit value (15) even though the
# wrapped species-name text changed length.
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 15
